$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.012.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.110.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.87%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.101.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.89%  "

$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("E10").Value = "  +13.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.93%  "

$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.624.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.932.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.109.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.96%  "

$ws.Range("E22").Value = "  +0.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.26%  "

$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.41%  "

$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("E29").Value = "  +5.03%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.37%  "

$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0860"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.10%  "

$ws.Range("E35").Value = "  +15.43%  "

$ws.Range("E36").Value = "  +4.50%  "

$ws.Range("E39").Value = "  +4.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "430.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.923.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0370"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.57%  "

$ws.Range("E45").Value = "  +5.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.08%  "

$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.00%  "

# Row 37/38: Filecoin and dogwifhat swap positions (with updated price/volume)
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +18.66%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.26%  "
